$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5. This shifts rows 5-9 down to 6-10.
# The previously-empty row 7 (B7/C7 blank placeholders) lands on row 8,
# which is exactly where the second milestone header belongs.
$ws.Rows("5:5").Insert() | Out-Null

# Update the activity codes to reference their milestone ("M - A" -> "M1 - A", etc.)
$ws.Range("C6").Value = "M1 - A"
$ws.Range("C7").Value = "M1 - B"
$ws.Range("C9").Value = "M2 - C"
$ws.Range("C10").Value = "M2 - D"

# Fill in the new milestone header rows.
$ws.Range("B5").Value = "Task 1"
$ws.Range("C5").Value = "M1"
$ws.Range("B8").Value = "Task 2"
$ws.Range("C8").Value = "M2"

# Style the new milestone header cells like the other header cells (bold white
# text on the blue fill used elsewhere in the sheet) plus word-wrap. Build the
# format once on B5, then fan it out to the other three cells so they all
# share a single new cell style, and merge each single cell (mirrors the
# merge entries the original authoring tool added).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").HorizontalAlignment = 1
$ws.Range("B5").WrapText = $true

foreach ($addr in @("C5", "B8", "C8")) {
    $ws.Range("B5").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

foreach ($addr in @("B5", "C5", "B8", "C8")) {
    $ws.Range($addr).MergeCells = $true
}

$excel.CutCopyMode = 0
